# Apply the skills-score corrections on the "Funcionarios" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (func1755382364785 / VIctor Ribeiro): Área 10 -> 7, Paciência 4 -> 8
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 8

# Row 3 (func1755382401524 / banana): Área was stored as the text "3",
# replace it with the numeric value 10 (also drops the now-unused
# shared string "3" from the workbook's shared string table).
$ws.Range("C3").Value = 10
